# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Rule row 11 ("R40") now resolves its computed/returned value, so the
# rule-name cell is updated from the literal label "R40" to "1".
# Keep it stored as text (matches the other rule-id cells in column B,
# e.g. "R10"/"R20"/"R30") rather than letting Excel auto-convert the
# digit string into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
